# Edit the quiz marksheet: recompute stats with corrected (numeric) marking
# scheme, drop the third (G/H) answer block and the now-unused D/E rows
# (19-40), and mark which answers were correct by echoing the correct
# answer into the "Student Ans" column for matching rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# ---- Summary rows (10-12): give the row-label cells in column A the
# same "mtitleStyle" formatting already used by the header row (A9),
# fix up the Right/Wrong/NotAttempt/Max counts and the Total/fraction. ----
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "No."

$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Marking"

$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Total"

$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 80
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = "80/112"

# ---- Third answer block (columns G/H) is no longer used - drop it. ----
$ws.Range("G15:H21").Clear()

# ---- Second answer block now only has 3 questions (rows 16-18); the
# remainder of its rows (19-40) are no longer used - drop them. ----
$ws.Range("D19:E40").Clear()

# ---- First & second blocks: where the student's answer matches the
# correct answer, echo the correct answer into the "Student Ans" cell
# (formatted with "correctStyle", same as B10/E16 etc.) ----
$ws.Range("B10").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Option A"

$ws.Range("B10").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = "Option D"

$ws.Range("B10").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Option B"

$ws.Range("B10").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = "Option C"

$ws.Range("B10").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Option C"

$ws.Range("B10").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Option D"

$ws.Range("B10").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A24").Value = "Option A"

$ws.Range("B10").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A25").Value = "Option A"

$ws.Range("B10").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A27").Value = "Option A"

$ws.Range("B10").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A28").Value = "Option D"

$ws.Range("B10").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = "Option D"

$ws.Range("B10").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A30").Value = "Option B"

$ws.Range("B10").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = "Option C"

$ws.Range("B10").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A33").Value = "Option D"

$ws.Range("B10").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A35").Value = "Option D"

$ws.Range("B10").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = "Option A"

$ws.Range("B10").Copy()
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("A39").Value = "Option D"

$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"

$ws.Range("B10").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "Option C"

$ws.Range("B10").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Option D"
